# BBox transforms in detection file.
# Adds two new rows ("cache_rate", "ds_type") to the dataset_params sheet,
# just above the existing "plan" row (which shifts from row 13 to row 15).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("dataset_params")
$ws.Activate()

# Insert two new blank rows above the current row 13 ("plan" row), pushing
# "plan" (and everything below it) down to row 15.
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# Match the row height used throughout the rest of the sheet.
$ws.Rows.Item(13).RowHeight = 13.8
$ws.Rows.Item(14).RowHeight = 13.8

# New row 13: cache_rate = 0.3
$ws.Range("A13").Value = "cache_rate"
$ws.Range("B13").Value = 0.3

# New row 14: ds_type = None
$ws.Range("A14").Value = "ds_type"
$ws.Range("B14").Value = "None"

# Restore the selection to where it now points after the insert (E15, which
# used to be E13 before the two new rows were added).
$ws.Range("E15").Select()

# Refresh the remembered selection on the other sheets too (cosmetic, no
# data change - just re-selecting the same cell that was already selected).
$wsModel = $wb.Worksheets.Item("model_params")
$wsModel.Activate()
$wsModel.Range("E7").Select()

$wsTransform = $wb.Worksheets.Item("transform_factors")
$wsTransform.Activate()
$wsTransform.Range("G7").Select()

$wsAffine = $wb.Worksheets.Item("affine3d")
$wsAffine.Activate()
$wsAffine.Range("E7").Select()

$wsLoss = $wb.Worksheets.Item("loss_params")
$wsLoss.Activate()
$wsLoss.Range("A7").Select()

$wsPlan1 = $wb.Worksheets.Item("plan1")
$wsPlan1.Activate()
$wsPlan1.Range("C48").Select()

$wsPlan2 = $wb.Worksheets.Item("plan2")
$wsPlan2.Activate()
$wsPlan2.Range("B1").Select()

# dataset_params was (and remains) the active/selected tab.
$ws.Activate()
